$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: female reactionTime
$ws.Range("B2").Value = "reactionTime"
$ws.Range("C2").Value = 43.70860258220497
$ws.Range("D2").Value = [double]"3.2269894340711e-10"
$ws.Range("E2").Value = $true
$ws.Range("F2").Value = [double]"3.220414182037056e-10"
$ws.Range("G2").Value = "Dunn"
$ws.Range("H2").Value = "April16"
$ws.Range("I2").Value = "June26"
$ws.Range("J2").Value = -0.7112482853223594
$ws.Range("K2").Value = 1.430864197530866

# Row 3: female peakTime
$ws.Range("B3").Value = "peakTime"
$ws.Range("C3").Value = 6.602513912876694
$ws.Range("D3").Value = 0.03683683599089073
$ws.Range("E3").Value = $true
$ws.Range("F3").Value = 0.049822393026943
$ws.Range("G3").Value = "Dunn"
$ws.Range("H3").Value = "June26"
$ws.Range("I3").Value = "May20"
$ws.Range("J3").Value = 0.2301097393689986
$ws.Range("K3").Value = -0.6074074074074005

# Row 4: female difference
$ws.Range("B4").Value = "difference"
$ws.Range("C4").Value = 41.64456638097869
$ws.Range("D4").Value = [double]"9.05726990314957e-10"
$ws.Range("E4").Value = $true
$ws.Range("F4").Value = [double]"4.913809571461353e-08"
$ws.Range("G4").Value = "Dunn"
$ws.Range("H4").Value = "April16"
$ws.Range("I4").Value = "June26"
$ws.Range("J4").Value = 0.5857338820301783
$ws.Range("K4").Value = -2.019753086419751

# Row 5: female peakValue
$ws.Range("B5").Value = "peakValue"
$ws.Range("C5").Value = 70.00295046057281
$ws.Range("D5").Value = [double]"6.295822118497218e-16"
$ws.Range("E5").Value = $true
$ws.Range("F5").Value = [double]"2.517974195281362e-15"
$ws.Range("G5").Value = "Dunn"
$ws.Range("H5").Value = "April16"
$ws.Range("I5").Value = "June26"
$ws.Range("J5").Value = 0.8789437585733882
$ws.Range("K5").Value = -43.05123456790125

# Row 6: female RMS
$ws.Range("B6").Value = "RMS"
$ws.Range("C6").Value = 78.10844503238788
$ws.Range("D6").Value = [double]"1.093872404052822e-17"
$ws.Range("E6").Value = $true
$ws.Range("F6").Value = [double]"1.811058574809449e-17"
$ws.Range("G6").Value = "Dunn"
$ws.Range("H6").Value = "April16"
$ws.Range("I6").Value = "June26"
$ws.Range("J6").Value = 0.8936899862825789
$ws.Range("K6").Value = -13.43082407407408

# Row 7: female tau (not significant, posthoc cells blank)
$ws.Range("B7").Value = "tau"
$ws.Range("C7").Value = 3.399749215245777
$ws.Range("D7").Value = 0.182706432610324
$ws.Range("E7").Value = $false
$ws.Range("F7:K7").ClearContents()

# Row 8: female AUC
$ws.Range("B8").Value = "AUC"
$ws.Range("C8").Value = 92.45678170786096
$ws.Range("D8").Value = [double]"8.380403795423334e-21"
$ws.Range("E8").Value = $true
$ws.Range("F8").Value = [double]"8.009213072298238e-21"
$ws.Range("G8").Value = "Dunn"
$ws.Range("H8").Value = "April16"
$ws.Range("I8").Value = "June26"
$ws.Range("J8").Value = 0.9403292181069959
$ws.Range("K8").Value = -6202.252033950618

# Row 9: male reactionTime
$ws.Range("B9").Value = "reactionTime"
$ws.Range("C9").Value = 6.04953810728265
$ws.Range("D9").Value = 0.0485690368170625
$ws.Range("E9").Value = $true
$ws.Range("F9").Value = 0.04532128050298975
$ws.Range("G9").Value = "Dunn"
$ws.Range("H9").Value = "April16"
$ws.Range("I9").Value = "May20"
$ws.Range("J9").Value = -0.2479423868312758
$ws.Range("K9").Value = 0.2191358024691343

# Row 10: male peakTime
$ws.Range("B10").Value = "peakTime"
$ws.Range("C10").Value = 7.310487494641959
$ws.Range("D10").Value = 0.02585519458062661
$ws.Range("E10").Value = $true
$ws.Range("F10").Value = 0.05618623759864159
$ws.Range("G10").Value = "Dunn"
$ws.Range("H10").Value = "April16"
$ws.Range("I10").Value = "June26"
$ws.Range("J10").Value = -0.2235939643347051
$ws.Range("K10").Value = 0.9790123456790134

# Row 11: male difference (not significant, posthoc cells blank)
$ws.Range("B11").Value = "difference"
$ws.Range("C11").Value = 1.652426416423395
$ws.Range("D11").Value = 0.4377036445564336
$ws.Range("E11").Value = $false
$ws.Range("F11:K11").ClearContents()

# Row 12: male peakValue
$ws.Range("B12").Value = "peakValue"
$ws.Range("C12").Value = 34.52994337542548
$ws.Range("D12").Value = [double]"3.176274776712328e-08"
$ws.Range("E12").Value = $true
$ws.Range("F12").Value = [double]"4.553856838281053e-08"
$ws.Range("G12").Value = "Dunn"
$ws.Range("H12").Value = "April16"
$ws.Range("I12").Value = "June26"
$ws.Range("J12").Value = 0.4670781893004116
$ws.Range("K12").Value = -64.82808641975306

# Row 13: male RMS
$ws.Range("B13").Value = "RMS"
$ws.Range("C13").Value = 15.21319692731006
$ws.Range("D13").Value = 0.0004971600936808299
$ws.Range("E13").Value = $true
$ws.Range("F13").Value = 0.0002966025195189699
$ws.Range("G13").Value = "Dunn"
$ws.Range("H13").Value = "April16"
$ws.Range("I13").Value = "June26"
$ws.Range("J13").Value = 0.4128943758573388
$ws.Range("K13").Value = -14.39813580246913

# Row 14: male tau
$ws.Range("B14").Value = "tau"
$ws.Range("C14").Value = 13.64772316056116
$ws.Range("D14").Value = 0.001087513283590059
$ws.Range("E14").Value = $true
$ws.Range("F14").Value = 0.0008544398307826692
$ws.Range("G14").Value = "Dunn"
$ws.Range("H14").Value = "April16"
$ws.Range("I14").Value = "May20"
$ws.Range("J14").Value = -0.4039780521262002
$ws.Range("K14").Value = 28.49903395061725

# Row 15: male AUC
$ws.Range("B15").Value = "AUC"
$ws.Range("C15").Value = 8.068671261581926
$ws.Range("D15").Value = 0.01769743383846177
$ws.Range("E15").Value = $true
$ws.Range("F15").Value = 0.01685596818818796
$ws.Range("G15").Value = "Dunn"
$ws.Range("H15").Value = "April16"
$ws.Range("I15").Value = "June26"
$ws.Range("J15").Value = 0.3161865569272977
$ws.Range("K15").Value = -5053.158296296298
